$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the price / volume columns with the latest scraped figures. The
# values are plain text in the sheet (e.g. "-0.68%", "300.36") rather than
# numeric/percentage cells, so force a text NumberFormat on each target cell
# right before writing it -- this keeps Excel from reinterpreting the text
# as a number/percentage and rewriting it in scientific/decimal form.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "300.36"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.68%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "38.14"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "8.31%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.014"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.57%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07718"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.63%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.185"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-6.46%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.963"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.94%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.997"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.34%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9174"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.62%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09055"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-8.78%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1795"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.07%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08452"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.77%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03554"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "7.05%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09937"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.20%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001477"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.45%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005665"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.25%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.477"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.216"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.30%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1319"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.19%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.568"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "6.18%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2251"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.25%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04663"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.31%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001232"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.29%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004431"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.22%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.01%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004763"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "40.30%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01740"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-2.15%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04682"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.35%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007923"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.81%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1385"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.94%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007701"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "11.91%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002292"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "10.18%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009760"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.28%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006060"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.89%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.17%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.764"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "194.73%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "34.78%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002106"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.17%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.17%"
